$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '59.683.79'
$ws.Range("E2").Value = '  +8.43%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.579.42'
$ws.Range("E3").Value = '  +10.22%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  -0.06%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '505.38'
$ws.Range("E5").Value = '  +6.31%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '156.84'
$ws.Range("E6").Value = '  +8.07%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.996'
$ws.Range("E7").Value = '  -0.28%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.610'
$ws.Range("E8").Value = '  -1.80%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '2.578.05'
$ws.Range("E9").Value = '  +10.08%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '6.09'
$ws.Range("E10").Value = '  +12.32%  '
$ws.Range("E11").Value = '  +6.90%  '
$ws.Range("E12").Value = '  +5.60%  '
$ws.Range("E13").Value = '  +1.46%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '2.996.81'
$ws.Range("E14").Value = '  +9.08%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '59.388.68'
$ws.Range("E15").Value = '  +7.92%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '21.82'
$ws.Range("E16").Value = '  +9.39%  '
$ws.Range("E17").Value = '  +5.79%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.574.04'
$ws.Range("E18").Value = '  +10.00%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '4.76'
$ws.Range("E19").Value = '  +4.40%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '339.86'
$ws.Range("E20").Value = '  +8.11%  '
$ws.Range("E21").Value = '  +7.95%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.04'
$ws.Range("E22").Value = '  +7.68%  '
$ws.Range("E23").Value = '  +0.02%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '59.99'
$ws.Range("E24").Value = '  +5.73%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.417'
$ws.Range("E25").Value = '  +6.15%  '
$ws.Range("E26").Value = '  +9.10%  '
$ws.Range("B27").Value = 'WrappedeETH'
$ws.Range("C27").Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.684.88'
$ws.Range("E27").Value = '  +9.86%  '
$ws.Range("B28").Value = 'Binance-PegBSC-USD'
$ws.Range("C28").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.998'
$ws.Range("E28").Value = '  -0.29%  '
$ws.Range("E29").Value = '  +11.96%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '7.34'
$ws.Range("E30").Value = '  +4.54%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.997'
$ws.Range("E31").Value = '  -0.28%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '157.13'
$ws.Range("E32").Value = '  +8.84%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '19.44'
$ws.Range("E33").Value = '  +7.04%  '
$ws.Range("E34").Value = '  +6.89%  '
$ws.Range("E35").Value = '  +8.17%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.20'
$ws.Range("E36").Value = '  +10.25%  '
$ws.Range("E37").Value = '  +9.29%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.858'
$ws.Range("E38").Value = '  +5.71%  '
$ws.Range("B39").Value = 'Bittensor'
$ws.Range("C39").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '301.66'
$ws.Range("E39").Value = '  +20.81%  '
$ws.Range("B40").Value = 'Filecoin'
$ws.Range("C40").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.73'
$ws.Range("E40").Value = '  +9.83%  '
$ws.Range("E41").Value = '  +9.18%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '35.22'
$ws.Range("E42").Value = '  +4.68%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.0576'
$ws.Range("E43").Value = '  +11.23%  '
$ws.Range("B44").Value = 'Mantle'
$ws.Range("C44").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.630'
$ws.Range("E44").Value = '  +9.80%  '
$ws.Range("B45").Value = 'Stellar'
$ws.Range("C45").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.102'
$ws.Range("E45").Value = '  +1.82%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.786'
$ws.Range("E46").Value = '  +25.45%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.00'
$ws.Range("E47").Value = '  +0.28%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '4.93'
$ws.Range("E48").Value = '  +13.65%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '19.22'
$ws.Range("E49").Value = '  +15.45%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0236'
$ws.Range("E50").Value = '  +7.95%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '10.26'
$ws.Range("E51").Value = '  +1.03%  '
